$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.001.50"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.52"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.64"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5092"
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3813"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08225"
$ws.Range("E9").Value = "  -7.77%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.59"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.188"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.44"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.851.23"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.174"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.43"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06608"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.007"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.022.96"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.226"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.556"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.075.74"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.69"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.43"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.607"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.595"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06528"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02409"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2169"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.202"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.242"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6410"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.856"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6090"
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.01"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.977"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.57"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.76"
$ws.Range("E51").Value = "  +1.04%  "
